$wb = $excel.ActiveWorkbook

# Sheet ALC, row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 319
$ws.Range("I96").Value = 319
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 957
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 416
$ws.Range("N96").ClearContents()

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 41667504
$ws.Range("I107").Value = 166666770
$ws.Range("K107").Value = 166666770
$ws.Range("M107").Value = -166664850

# Sheet ALC, row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 14495724
$ws.Range("I125").Value = 1890.2858
$ws.Range("K125").Value = 17012.5722
$ws.Range("M125").Value = -14552.5722

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1960.5862
$ws.Range("I132").Value = 2085.5557
$ws.Range("J132").Value = 273.5
$ws.Range("K132").Value = 6256.6671
$ws.Range("L132").Value = 820.5
$ws.Range("M132").Value = -3726.6671
$ws.Range("N132").Value = -5880.5

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 39891.207
$ws.Range("I137").Value = 54479.676
$ws.Range("K137").Value = 163439.028
$ws.Range("M137").Value = -160889.028

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3629.5334
$ws.Range("J138").Value = 3947.4443
$ws.Range("L138").Value = 11842.3329
$ws.Range("N138").Value = -22122.3329

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2829546
$ws.Range("I2").Value = 5143451.5
$ws.Range("J2").Value = 1439.1111
$ws.Range("K2").Value = 5143451.5
$ws.Range("L2").Value = 1439.1111
$ws.Range("M2").Value = -5143338.5
$ws.Range("N2").Value = -1665.1111

# Sheet ARM, row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 189.14285
$ws.Range("I5").Value = 189.14285
$ws.Range("K5").Value = 189.14285
$ws.Range("M5").Value = -77.14285000000001

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17094.486
$ws.Range("I32").Value = 14551.692
$ws.Range("K32").Value = 14551.692
$ws.Range("M32").Value = -14264.692

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 40628.082
$ws.Range("J74").Value = 85041.63
$ws.Range("L74").Value = 85041.63
$ws.Range("N74").Value = -86789.63

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 40628.082
$ws.Range("J77").Value = 85041.63
$ws.Range("L77").Value = 425208.15
$ws.Range("N77").Value = -433944.15

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2829546
$ws.Range("I116").Value = 5143451.5
$ws.Range("J116").Value = 1439.1111
$ws.Range("K116").Value = 5143451.5
$ws.Range("L116").Value = 1439.1111
$ws.Range("M116").Value = -5141157.5
$ws.Range("N116").Value = -6027.1111

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1752146.1
$ws.Range("I122").Value = 1596785.6
$ws.Range("K122").Value = 4790356.800000001
$ws.Range("M122").Value = -4787906.800000001

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2829546
$ws.Range("I3").Value = 5143451.5
$ws.Range("J3").Value = 1439.1111
$ws.Range("K3").Value = 5143451.5
$ws.Range("L3").Value = 1439.1111
$ws.Range("M3").Value = -5143337.5
$ws.Range("N3").Value = -1667.1111

# Sheet BSM, row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 189.14285
$ws.Range("I4").Value = 189.14285
$ws.Range("K4").Value = 189.14285
$ws.Range("M4").Value = -74.14285000000001

# Sheet BSM, row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 9997
$ws.Range("J81").Value = 9997
$ws.Range("L81").Value = 9997
$ws.Range("N81").Value = -12119

# Sheet BSM, row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 9997
$ws.Range("J84").Value = 9997
$ws.Range("L84").Value = 29991
$ws.Range("N84").Value = -40599

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5094.9688
$ws.Range("I134").Value = 1913.4138
$ws.Range("K134").Value = 5740.2414
$ws.Range("M134").Value = -3205.2414

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 852.65
$ws.Range("J16").Value = 1336.25
$ws.Range("L16").Value = 1336.25
$ws.Range("N16").Value = -1910.25

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19611.885
$ws.Range("I31").Value = 3137.8064
$ws.Range("J31").Value = 36635.1
$ws.Range("K31").Value = 3137.8064
$ws.Range("L31").Value = 36635.1
$ws.Range("M31").Value = -2842.8064
$ws.Range("N31").Value = -37225.1

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 19611.885
$ws.Range("I34").Value = 3137.8064
$ws.Range("J34").Value = 36635.1
$ws.Range("K34").Value = 3137.8064
$ws.Range("L34").Value = 36635.1
$ws.Range("M34").Value = -2935.8064
$ws.Range("N34").Value = -37039.1

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9825.200000000001
$ws.Range("I58").Value = 13237.556
$ws.Range("J58").Value = 4706.6665
$ws.Range("K58").Value = 13237.556
$ws.Range("L58").Value = 4706.6665
$ws.Range("M58").Value = -13034.556
$ws.Range("N58").Value = -5112.6665

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5852.75
$ws.Range("I86").Value = 4712.567
$ws.Range("K86").Value = 4712.567
$ws.Range("M86").Value = -3589.567

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5852.75
$ws.Range("I89").Value = 4712.567
$ws.Range("K89").Value = 23562.835
$ws.Range("M89").Value = -17946.835

# Sheet CRP, row 93
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 13886
$ws.Range("I93").Value = 10840.4
$ws.Range("J93").Value = 21500
$ws.Range("K93").Value = 10840.4
$ws.Range("L93").Value = 21500
$ws.Range("M93").Value = -8968.4
$ws.Range("N93").Value = -25244

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 852.65
$ws.Range("J113").Value = 1336.25
$ws.Range("L113").Value = 1336.25
$ws.Range("N113").Value = -5676.25

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 9825.200000000001
$ws.Range("I136").Value = 13237.556
$ws.Range("J136").Value = 4706.6665
$ws.Range("K136").Value = 39712.66800000001
$ws.Range("L136").Value = 14119.9995
$ws.Range("M136").Value = -37162.66800000001
$ws.Range("N136").Value = -19219.9995

# Sheet CUL, row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 132283.89
$ws.Range("I46").Value = 479211.72
$ws.Range("J46").Value = 4468.3687
$ws.Range("K46").Value = 1437635.16
$ws.Range("L46").Value = 13405.1061
$ws.Range("M46").Value = -1437544.16
$ws.Range("N46").Value = -13587.1061

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1829.55
$ws.Range("I132").Value = 1708.0834
$ws.Range("J132").Value = 2011.75
$ws.Range("K132").Value = 15372.7506
$ws.Range("L132").Value = 18105.75
$ws.Range("M132").Value = -12842.7506
$ws.Range("N132").Value = -23165.75

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1746934.6
$ws.Range("J80").Value = 417083.16
$ws.Range("L80").Value = 417083.16
$ws.Range("N80").Value = -419079.16

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1746934.6
$ws.Range("J83").Value = 417083.16
$ws.Range("L83").Value = 2085415.8
$ws.Range("N83").Value = -2095399.8

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 37052436
$ws.Range("I93").Value = 66669984
$ws.Range("K93").Value = 66669984
$ws.Range("M93").Value = -66668736

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5880
$ws.Range("I122").Value = 4910.65
$ws.Range("K122").Value = 14731.95
$ws.Range("M122").Value = -12281.95

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 33334690
$ws.Range("I107").Value = 52632170
$ws.Range("K107").Value = 157896510
$ws.Range("M107").Value = -157894590

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3449.3784
$ws.Range("I122").Value = 1949.5555
$ws.Range("J122").Value = 7498.9
$ws.Range("K122").Value = 5848.666499999999
$ws.Range("L122").Value = 22496.7
$ws.Range("M122").Value = -3398.666499999999
$ws.Range("N122").Value = -27396.7

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17137040
$ws.Range("I132").Value = 17860562
$ws.Range("K132").Value = 53581686
$ws.Range("M132").Value = -53579156

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2988.22
$ws.Range("I136").Value = 2697.55
$ws.Range("K136").Value = 8092.650000000001
$ws.Range("M136").Value = -5542.650000000001
